$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column C
$ws.Range("C1").Value = "Value"

# Values for column C, rows 2-17
$values = @{
    2  = 2
    3  = 5
    4  = 10
    5  = 65
    6  = 27
    7  = 58
    8  = 108
    9  = 134
    10 = 209
    11 = 250
    12 = 280
    13 = 395
    14 = 674
    15 = 538
    16 = 1624
    17 = 2048
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}

# Copy style from B1:B17 into C1:C17 so style index matches (s="1")
$ws.Range("B1:B17").Copy()
$ws.Range("C1:C17").PasteSpecial(-4122) # xlPasteFormats

# Set selection to E5 as per diff
$ws.Range("E5").Select()

# Page setup: paper size A4(9) Letter? -> paperSize 9, orientation portrait
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
